$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "MathNet.Numerics"
$ws.Range("C10").Value = "MIT/X11"
$ws.Range("D10").Value = "https://numerics.mathdotnet.com/License.html"

$ws.Range("A10").Borders.Item(7).LineStyle = 1
$ws.Range("A10").Borders.Item(10).LineStyle = 1
$ws.Range("C10").Borders.Item(7).LineStyle = 1
$ws.Range("C10").Borders.Item(10).LineStyle = 1

$ws.Range("E16").Select()
